# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) for the crafting leves whose market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip (Eye Drops)
$ws.Range("H8").Value = 507.33334
$ws.Range("I8").Value = 50
$ws.Range("J8").Value = 964.6667
$ws.Range("K8").Value = 150
$ws.Range("L8").Value = 2894.0001
$ws.Range("M8").Value = -11
$ws.Range("N8").Value = -3172.0001

# Row 40: Stuck in the Moment (Horn Glue)
$ws.Range("H40").Value = 1624.85
$ws.Range("I40").Value = 1627.6111
$ws.Range("K40").Value = 1627.6111
$ws.Range("M40").Value = -1452.6111

# Row 41: The Write Stuff (Enchanted Mythril Ink)
$ws.Range("H41").Value = 2667.375
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# Row 70: Consecrating Congregation (Holy Water)
$ws.Range("H70").Value = 9688.667
$ws.Range("I70").Value = 10199.667
$ws.Range("J70").Value = 8666.667
$ws.Range("K70").Value = 30599.001
$ws.Range("L70").Value = 26000.001
$ws.Range("M70").Value = -30329.001
$ws.Range("N70").Value = -26540.001

# Row 73: Curbing the Contagion (L) (Holy Water)
$ws.Range("H73").Value = 9688.667
$ws.Range("I73").Value = 10199.667
$ws.Range("J73").Value = 8666.667
$ws.Range("K73").Value = 30599.001
$ws.Range("L73").Value = 26000.001
$ws.Range("M73").Value = -29663.001
$ws.Range("N73").Value = -27872.001

# Row 125: Body over Mind (Grade 5 Dexterity Alkahest)
$ws.Range("H125").Value = 3374.5
$ws.Range("I125").Value = 3374.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 30370.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -27910.5
$ws.Range("N125").ClearContents()

# Row 132: Fast-forwarding Flora (Growth Formula Lambda)
$ws.Range("H132").Value = 1735.6316
$ws.Range("I132").Value = 1776.6666
$ws.Range("J132").Value = 997
$ws.Range("K132").Value = 5329.9998
$ws.Range("L132").Value = 2991
$ws.Range("M132").Value = -2799.9998
$ws.Range("N132").Value = -8051

# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 919.75
$ws.Range("I137").Value = 919.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2759.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -209.25
$ws.Range("N137").ClearContents()

# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 3308.5833
$ws.Range("I138").Value = 1965.8182
$ws.Range("J138").Value = 3707.7837
$ws.Range("K138").Value = 5897.4546
$ws.Range("L138").Value = 11123.3511
$ws.Range("M138").Value = -757.4546
$ws.Range("N138").Value = -21403.3511

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 2657750.5
$ws.Range("I32").Value = 2595324.5
$ws.Range("K32").Value = 2595324.5
$ws.Range("M32").Value = -2595037.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 858.2222
$ws.Range("J31").Value = 832.5
$ws.Range("L31").Value = 832.5
$ws.Range("N31").Value = -1422.5

# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 858.2222
$ws.Range("J34").Value = 832.5
$ws.Range("L34").Value = 832.5
$ws.Range("N34").Value = -1236.5

# Row 42: Live Freelance or Die (Heavy Steel Lance)
$ws.Range("H42").Value = 37500
$ws.Range("I42").Value = 37500
$ws.Range("J42").Value = 37500
$ws.Range("K42").Value = 37500
$ws.Range("L42").Value = 37500
$ws.Range("M42").Value = -36907
$ws.Range("N42").Value = -38686

# Row 50: The Arsenal of Theocracy (Cobalt Halberd)
$ws.Range("H50").Value = 22510.375
$ws.Range("I50").Value = 5083
$ws.Range("K50").Value = 5083
$ws.Range("M50").Value = -4458

# Row 58: You Do the Heavy Lifting (Mahogany Lumber)
$ws.Range("H58").Value = 1811.7894
$ws.Range("I58").Value = 1127.1538
$ws.Range("K58").Value = 1127.1538
$ws.Range("M58").Value = -924.1538

# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 6663.5713
$ws.Range("I132").Value = 6663.5713
$ws.Range("K132").Value = 19990.7139
$ws.Range("M132").Value = -17460.7139

# Row 134: Wood You Be Quiet (Ceiba Lumber)
$ws.Range("H134").Value = 1292.2858
$ws.Range("I134").Value = 1276.3846
$ws.Range("K134").Value = 3829.1538
$ws.Range("M134").Value = -1294.1538

# Row 136: Turali Quality (Dark Mahogany Lumber)
$ws.Range("H136").Value = 1811.7894
$ws.Range("I136").Value = 1127.1538
$ws.Range("K136").Value = 3381.4614
$ws.Range("M136").Value = -831.4614000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch (Chamomile Tea)
$ws.Range("H34").Value = 2326.9092
$ws.Range("J34").Value = 2221.7778
$ws.Range("L34").Value = 6665.3334
$ws.Range("N34").Value = -6833.3334

# Row 39: Bloody Good Tart, This (Blood Currant Tart)
$ws.Range("H39").Value = 1630.4166
$ws.Range("I39").Value = 891.25
$ws.Range("K39").Value = 2673.75
$ws.Range("M39").Value = -2379.75

$ws = $wb.Worksheets.Item("GSM")
# Row 26: Perk of Fiction (Coral Ring)
$ws.Range("H26").Value = 48974.5
$ws.Range("J26").Value = 48974.5
$ws.Range("L26").Value = 48974.5
$ws.Range("N26").Value = -49534.5

# Row 50: Coral on My Mind (Red Coral Ring)
$ws.Range("H50").Value = 48974.5
$ws.Range("J50").Value = 48974.5
$ws.Range("L50").Value = 48974.5
$ws.Range("N50").Value = -49970.5

# Row 70: Sky Is the Limit (Mythrite Ingot)
$ws.Range("H70").Value = 5444
$ws.Range("I70").Value = 5444
$ws.Range("K70").Value = 5444
$ws.Range("M70").Value = -5174

# Row 73: Hulls of Broken Dreams (L) (Mythrite Ingot)
$ws.Range("H73").Value = 5444
$ws.Range("I73").Value = 5444
$ws.Range("K73").Value = 5444
$ws.Range("M73").Value = -4508

$ws = $wb.Worksheets.Item("LTW")
# Row 56: Hold On Tight (Raptorskin Smithy's Gloves)
$ws.Range("H56").Value = 4175017
$ws.Range("I56").Value = 4175017
$ws.Range("K56").Value = 4175017
$ws.Range("M56").Value = -4174326

# Row 61: Spelling Me Softly (Raptor Leather)
$ws.Range("H61").Value = 1032.7273
$ws.Range("I61").Value = 929.1111
$ws.Range("K61").Value = 929.1111
$ws.Range("M61").Value = -727.1111

# Row 113: Peace in Rest (Atrociraptor Leather)
$ws.Range("H113").Value = 1032.7273
$ws.Range("I113").Value = 929.1111
$ws.Range("K113").Value = 929.1111
$ws.Range("M113").Value = 1240.8889

$ws = $wb.Worksheets.Item("WVR")
# Row 46: Crunching the Numbers (Linen Hat)
$ws.Range("H46").Value = 99429
$ws.Range("J46").Value = 99429
$ws.Range("L46").Value = 99429
$ws.Range("N46").Value = -99891

# Row 61: Bundle Up, It's Odd out There (Woolen Deerstalker)
$ws.Range("H61").Value = 28016.666
$ws.Range("J61").Value = 44999
$ws.Range("L61").Value = 44999
$ws.Range("N61").Value = -45583

# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 2462.6858
$ws.Range("I122").Value = 1710.0476
$ws.Range("J122").Value = 3591.6428
$ws.Range("K122").Value = 5130.142800000001
$ws.Range("L122").Value = 10774.9284
$ws.Range("M122").Value = -2680.142800000001
$ws.Range("N122").Value = -15674.9284

# Row 123: Helping Handwear (Fingerless Darkhempen Gloves of Healing)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 134: Cloth for Canvas (Mountain Linen)
$ws.Range("H134").Value = 99429
$ws.Range("J134").Value = 99429
$ws.Range("L134").Value = 298287
$ws.Range("N134").Value = -303357

# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 2470
$ws.Range("I136").Value = 2364.625
$ws.Range("K136").Value = 7093.875
$ws.Range("M136").Value = -4543.875
